# Airbase Input Specification Document - robustness update
# - Swap the "Can Use?"/"Reason:" info between the Forward Slash (row 15)
#   and Backward Slash (row 16) rows so the escaping caveat is attached to
#   the Backward Slash character instead of the Forward Slash character.
# - Reword the "Special Characters" reason text (row 19) to state they are
#   generally NOT accepted by MySQL.
# - Update the view/selection state to reflect where the user left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 15 / Row 16: swap Can Use? (col C) and Reason (col D) ---
# Before: Row15 (Forward Slash)  -> C15=NO,  D15="Escape character in python and disagrees with MySQL"
#         Row16 (Backward Slash) -> C16=YES, D16=(blank)
# After:  Row15 (Forward Slash)  -> C15=YES, D15=(blank)
#         Row16 (Backward Slash) -> C16=NO,  D16="Escape character in python and disagrees with MySQL"
$ws.Range("C15").Value2 = "YES"
$ws.Range("D15").ClearContents()

$ws.Range("C16").Value2 = "NO"
$ws.Range("D16").Value2 = "Escape character in python and disagrees with MySQL"

# --- Row 19: Special Characters reason text update ---
$ws.Range("D19").Value2 = "Generally not accepted by MySQL"

# --- Sheet view / selection state ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D19").Select()
